# Regenerate save_data column G ("K", i.e. strikeouts) with freshly
# calculated values (std/mean recomputed, s_vals recalculated) for rows 2-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0,1,0,1,0,1,1,2,0,2,0,0,2,1,0,1,2,0,0,1,1,1,0,2,0,2,2,1,1,1,0,1,1,1,1,0,0,2,1,0,0,0,1,0,1,0,1,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
